$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sizes")

# Clear existing content/strings so stale shared strings are dropped on save.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Number of sequences"
$ws.Range("B1").Value = "Input FASTQ"
$ws.Range("C1").Value = "Pip Database"

# Data rows: number of sequences, input FASTQ size (MB), pip database size (MB)
$ws.Range("A2").Value = 500000
$ws.Range("B2").Value = 371.4
$ws.Range("C2").Value = 206.9

$ws.Range("A3").Value = 1000000
$ws.Range("B3").Value = 743
$ws.Range("C3").Value = 413.9

$ws.Range("A4").Value = 2000000
$ws.Range("B4").Value = 1486
$ws.Range("C4").Value = 828

$ws.Range("A5").Value = 4000000
$ws.Range("B5").Value = 3051.52
$ws.Range("C5").Value = 1699.84

$ws.Range("A6").Value = 8000000
$ws.Range("B6").Value = 6082.56
$ws.Range("C6").Value = 3389.44

$ws.Range("A7").Value = 16000000
$ws.Range("B7").Value = 12165.12
$ws.Range("C7").Value = 6789.12

$ws.Range("A8").Value = 32000000
$ws.Range("B8").Value = 24350.72

$ws.Range("A9").Value = 64000000
$ws.Range("B9").Value = 48701.44

$ws.Range("A10").Value = 91000000
$ws.Range("B10").Value = 69754.88
$ws.Range("C10").Value = 35061.76

# Column widths
$ws.Range("B:B").ColumnWidth = 31
$ws.Range("C:C").ColumnWidth = 28.33203125

# Selection
$null = $ws.Range("C2").Select()

Write-Host "data done"
